# repull data, push all data, mean calculation
# Updates column F ("dSF") values for a set of rows to reflect the
# refreshed/repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 2
    9  = 1
    10 = -2
    12 = -1
    15 = 0
    24 = -5
    26 = 1
    29 = 3
    34 = 3
    36 = 3
    37 = 0
    40 = 1
    44 = 0
    47 = -2
    48 = 5
    55 = -1
    60 = -3
    62 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
